# Realestate Update resale numbers 2023-06-12 22:19
# Append a new data row (42) to the CityResaleNum sheet with the latest
# resale-number snapshot, matching the existing row layout:
# A=Date, B=Time, C=Weekday, D=Week (text), E..T = numeric city values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 42

# Columns A:D hold plain text values in this sheet (e.g. "2023-06-12",
# "24" for the week number). Excel's COM layer auto-detects dates/numbers
# from plain strings, so force a text format while assigning them, then
# clear the formatting back to the sheet's default (no explicit style),
# just like the other data rows.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-12"
$ws.Cells.Item($row, 2).Value = "22:12:49"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "24"

$textRange.ClearFormats()

# Columns E:T are numeric city resale-number values.
$ws.Cells.Item($row, 5).Value = 121428
$ws.Cells.Item($row, 6).Value = 134969
$ws.Cells.Item($row, 7).Value = 161258
$ws.Cells.Item($row, 8).Value = 131947
$ws.Cells.Item($row, 9).Value = 176388
$ws.Cells.Item($row, 10).Value = 113928
$ws.Cells.Item($row, 11).Value = 202035
$ws.Cells.Item($row, 12).Value = 222839
$ws.Cells.Item($row, 13).Value = 174001
$ws.Cells.Item($row, 14).Value = 100797
$ws.Cells.Item($row, 15).Value = 38833
$ws.Cells.Item($row, 16).Value = 34101
$ws.Cells.Item($row, 17).Value = 51323
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36900
$ws.Cells.Item($row, 20).Value = -1
